$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data record was inserted at row 57, pushing every
# subsequent record (old rows 57..158) down by one (new rows 58..159).
# Insert a blank row at position 57 (shifts rows 57-158 down to 58-159,
# carrying their original data/styles with them automatically).
$ws.Rows.Item(57).EntireRow.Insert()

# Fill the newly inserted row 57 with the new record. It mirrors the
# data that used to be in row 57 (same market/product/variety/etc.)
# except for a new date (2021-12-03, serial 44533).
$ws.Range("A57").Value2 = 1
$ws.Range("B57").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C57").Value2 = "Arica y Parinacota"
$ws.Range("D57").Value2 = 44533
$ws.Range("E57").Value2 = 15
$ws.Range("F57").Value2 = "Fruta"
$ws.Range("G57").Value2 = 100108
$ws.Range("H57").Value2 = "Tropicales y subtropicales"
$ws.Range("I57").Value2 = 100108006
$ws.Range("J57").Value2 = "Plátano"
$ws.Range("K57").Value2 = "Sin especificar"
$ws.Range("L57").Value2 = "Pintón"
$ws.Range("M57").Value2 = 120
$ws.Range("N57").Value2 = 15000
$ws.Range("O57").Value2 = 16000
$ws.Range("P57").Value2 = 15500
$ws.Range("Q57").Value2 = "$/caja 20 kilos"
$ws.Range("R57").Value2 = "Ecuador"
$ws.Range("S57").Value2 = 775
$ws.Range("T57").Value2 = 20
